# fix: alterar python version para 3.11.5
# Update absenteeism data rows 2-11 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 63316
$ws.Range("B2").Value = "Sra. Valentina Pires"
$ws.Range("C2").Value = "TI"
$ws.Range("D2").Value = "Doença"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 45095
$ws.Range("G2").Value = 5812.29

# Row 3
$ws.Range("A3").Value = 195
$ws.Range("B3").Value = "Caio Farias"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 45106
$ws.Range("G3").Value = 6407.42

# Row 4
$ws.Range("A4").Value = 36616
$ws.Range("B4").Value = "Vicente da Mata"
$ws.Range("C4").Value = "Vendas"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 45082
$ws.Range("G4").Value = 6950.06

# Row 5
$ws.Range("A5").Value = 73965
$ws.Range("B5").Value = "Larissa Barbosa"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45106
$ws.Range("G5").Value = 5745.93

# Row 6
$ws.Range("A6").Value = 48273
$ws.Range("B6").Value = "Ana Julia Gomes"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Problemas pessoais"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45101
$ws.Range("G6").Value = 10747.78

# Row 7
$ws.Range("A7").Value = 65771
$ws.Range("B7").Value = "Rafaela Ribeiro"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Consulta médica"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45080
$ws.Range("G7").Value = 3872.27

# Row 8
$ws.Range("A8").Value = 56353
$ws.Range("B8").Value = "Bianca Nascimento"
$ws.Range("C8").Value = "TI"
$ws.Range("D8").Value = "Consulta médica"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45081
$ws.Range("G8").Value = 7468.14

# Row 9
$ws.Range("A9").Value = 29866
$ws.Range("B9").Value = "Emanuella da Cruz"
$ws.Range("C9").Value = "Vendas"
$ws.Range("D9").Value = "Viagem de negócios"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45093
$ws.Range("G9").Value = 3270.91

# Row 10
$ws.Range("A10").Value = 1860
$ws.Range("B10").Value = "Renan Pires"
$ws.Range("C10").Value = "Jurídico"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45098
$ws.Range("G10").Value = 2754.81

# Row 11
$ws.Range("A11").Value = 69796
$ws.Range("B11").Value = "Joaquim Cardoso"
$ws.Range("C11").Value = "TI"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 45100
$ws.Range("G11").Value = 6756.5
